# Update the "courses" worksheet: recategorise the department column (C),
# split the combined location/availability text (M) into location (M) and
# locationDetail (N) for the massage + massage-package rows, and clear the
# now-unused promotionValidity column (R).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# --- Column C: department -> course-category name ------------------------
$ws.Range("C2:C7").Value   = "Ageing Support"
$ws.Range("C8").Value      = "Community Services"
$ws.Range("C9:C10").Value  = "Early Childhood"
$ws.Range("C11:C12").Value = "Massage"
$ws.Range("C13:C20").Value = "Packages"

# --- Columns M/N: split "NSW/QLD/TAS (Currently not accepting enrolments)"
#     into location = "NSW/QLD/TAS" and locationDetail = "Currently not
#     accepting enrolments" for the Massage rows and the Massage package row.
foreach ($r in 11, 12, 20) {
    $ws.Cells.Item($r, 13).Value = "NSW/QLD/TAS"
    $ws.Cells.Item($r, 14).Value = "Currently not accepting enrolments"
}

# --- Column R: promotionValidity no longer applies - clear all values ----
$ws.Range("R2:R20").ClearContents()

$wb.Save()
